# Maestro.xlsx — "Add files via upload"
#
# A new article (cigarette product, barcode 7798100200538) is inserted as
# the first data row (row 2) of the "Artículos" sheet, pushing every
# existing article down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Insert a new blank row above the current row 2 (shifts rows 2..64 -> 3..65)
$ws.Rows.Item(2).Insert()

# Populate the new row with the new article's data
$ws.Range("A2").Value = 7798100200538
$ws.Range("B2").Value = "Cigarrillos"
$ws.Range("C2").Value = "rubios"
$ws.Range("D2").Value = "pink wave mintz"
$ws.Range("E2").Value = "Milenio"
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = "und."
$ws.Range("H2").Value = "caja"
$ws.Range("I2").Value = "Cigarrillos"
$ws.Range("J2").Value = "Argentina"
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = $false
$ws.Range("M2").Value = $false
$ws.Range("N2").Value = "C:\EditaSoft\Imágenes de artículos\7798100200538.png"
$ws.Range("O2").Value = $true
$ws.Range("P2").Value = $true
